$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 234, shifting the old rows 234-237 down to 236-239
$ws.Rows.Item(234).Resize(2).Insert()

# New row 234 (Magnum, Peru, higher prices)
$ws.Cells.Item(234, 1).Value = 9
$ws.Cells.Item(234, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(234, 3).Value = "Metropolitana"
$ws.Cells.Item(234, 4).Value = 44448
$ws.Cells.Item(234, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112031
$ws.Cells.Item(234, 7).Value = "Poroto verde"
$ws.Cells.Item(234, 8).Value = "Magnum"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 34
$ws.Cells.Item(234, 11).Value = 37000
$ws.Cells.Item(234, 12).Value = 38000
$ws.Cells.Item(234, 13).Value = 37500
$ws.Cells.Item(234, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(234, 15).Value = "Perú"
$ws.Cells.Item(234, 16).Value = 1500
$ws.Cells.Item(234, 17).Value = 25
$ws.Cells.Item(234, 18).Value = "Hortaliza"

# New row 235 (Sin especificar, Peru, higher prices)
$ws.Cells.Item(235, 1).Value = 9
$ws.Cells.Item(235, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(235, 3).Value = "Metropolitana"
$ws.Cells.Item(235, 4).Value = 44448
$ws.Cells.Item(235, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(235, 5).Value = 13
$ws.Cells.Item(235, 6).Value = 100112031
$ws.Cells.Item(235, 7).Value = "Poroto verde"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 16
$ws.Cells.Item(235, 11).Value = 36000
$ws.Cells.Item(235, 12).Value = 38000
$ws.Cells.Item(235, 13).Value = 37000
$ws.Cells.Item(235, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(235, 15).Value = "Perú"
$ws.Cells.Item(235, 16).Value = 1480
$ws.Cells.Item(235, 17).Value = 25
$ws.Cells.Item(235, 18).Value = "Hortaliza"
